$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Total shipping fee" line above the Subtotal row ---
# This pushes the old rows 10-28 (Subtotal, HST, Total, Notes, ...) down to 11-29.
$ws.Rows("10:10").Insert()

# Carry over the plain-data-row formatting (border/font, no number format) to A10:F10
# from the row directly above it (row 9), matching the rest of the line-item rows.
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)

# G10 should use the currency-with-cents format ($#,##0.00), same as G4 (and not the
# $#,##0 format some of the other price cells currently use).
$ws.Range("G4").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New line item content
$ws.Range("E10").Value = "Total shipping fee"
$ws.Range("G10").Value = 19.96

# The Subtotal formula (now on row 11 after the insert) needs to widen to include
# the new shipping-fee row.
$ws.Range("G11").Formula = "=SUM(G4:G10)"

# Normalize the price column's number format for the other line items (G7:G9) so they
# match G4:G6 ($#,##0.00 instead of $#,##0).
$ws.Range("G4").Copy()
$ws.Range("G7:G9").PasteSpecial(-4122)
$ws.Range("G7").Value = 16.99
$ws.Range("G8").Value = 1.78
$ws.Range("G9").Value = 6.89
$ws.Application.CutCopyMode = $false

$ws.Range("H11").Select()
